$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('ITA_Italy', 118)
    ,('FIN_Finland', 118)
    ,('CAN_Canada', 118)
    ,('CHE_Switzerland', 118)
    ,('ISR_Israel', 118)
    ,('ISL_Iceland', 118)
    ,('IRL_Ireland', 118)
    ,('IND_India', 118)
    ,('COL_Colombia', 118)
    ,('CRI_Costa Rica', 118)
    ,('GTM_Guatemala', 118)
    ,('THA_Thailand', 118)
    ,('CYP_Cyprus', 118)
    ,('DEU_Germany', 118)
    ,('NOR_Norway', 118)
    ,('SWE_Sweden', 118)
    ,('DNK_Denmark', 118)
    ,('MUS_Mauritius', 118)
    ,('ECU_Ecuador', 118)
    ,('EGY_Egypt', 118)
    ,('ESP_Spain', 118)
    ,('MEX_Mexico', 118)
    ,('GBR_United Kingdom', 118)
    ,('NZL_New Zealand', 118)
    ,('JPN_Japan', 118)
    ,('FRA_France', 118)
    ,('BOL_Bolivia (Plurinational State of)', 118)
    ,('PER_Peru', 118)
    ,('ZAF_South Africa', 118)
    ,('NLD_Netherlands', 118)
    ,('VEN_Venezuela (Bolivarian Republic of)', 118)
    ,('ARG_Argentina', 118)
    ,('USA_United States', 118)
    ,('URY_Uruguay', 118)
    ,('AUS_Australia', 118)
    ,('AUT_Austria', 118)
    ,('PHL_Philippines', 118)
    ,('BEL_Belgium', 118)
    ,('BRA_Brazil', 118)
    ,('LUX_Luxembourg', 118)
    ,('MAR_Morocco', 118)
    ,('PRT_Portugal', 118)
    ,('LKA_Sri Lanka', 118)
    ,('TTO_Trinidad and Tobago', 118)
    ,('KEN_Kenya', 118)
    ,('TUR_Turkey', 118)
    ,('CHL_Chile', 136)
    ,('GRC_Greece', 142)
    ,('PRY_Paraguay', 142)
    ,('TWN_Taiwan', 142)
    ,('DOM_Dominican Republic', 142)
    ,('CHN_China', 166)
    ,('NGA_Nigeria', 184)
    ,('KOR_Republic of Korea', 190)
    ,('JAM_Jamaica', 190)
    ,('PAN_Panama', 208)
    ,('MLT_Malta', 214)
    ,('JOR_Jordan', 214)
    ,('HND_Honduras', 214)
    ,('ZWE_Zimbabwe', 214)
    ,('ZMB_Zambia', 238)
    ,('MYS_Malaysia', 238)
    ,('IRN_Iran (Islamic Republic of)', 238)
    ,('NIC_Nicaragua', 274)
    ,('IDN_Indonesia', 334)
    ,('BFA_Burkina Faso', 337)
    ,('GAB_Gabon', 358)
    ,('TZA_U.R. of Tanzania: Mainland', 358)
    ,('ROU_Romania', 358)
    ,('RWA_Rwanda', 358)
    ,('TUN_Tunisia', 358)
    ,('BRB_Barbados', 358)
    ,('CIV_Côte d''Ivoire', 358)
    ,('CMR_Cameroon', 358)
    ,('BWA_Botswana', 358)
    ,('NER_Niger', 358)
    ,('HKG_China, Hong Kong SAR', 358)
    ,('SEN_Senegal', 358)
    ,('NAM_Namibia', 358)
    ,('SGP_Singapore', 358)
    ,('MOZ_Mozambique', 358)
    ,('MRT_Mauritania', 436)
    ,('BEN_Benin', 436)
    ,('CAF_Central African Republic', 454)
    ,('LSO_Lesotho', 454)
    ,('BDI_Burundi', 454)
    ,('FJI_Fiji', 454)
    ,('TGO_Togo', 454)
    ,('SLE_Sierra Leone', 472)
    ,('SLV_El Salvador', 492)
    ,('UGA_Uganda', 492)
    ,('COD_D.R. of the Congo', 492)
    ,('PAK_Pakistan', 492)
    ,('ETH_Ethiopia', 522)
    ,('MWI_Malawi', 564)
    ,('GHA_Ghana', 582)
    ,('QAT_Qatar', 598)
    ,('KWT_Kuwait', 598)
    ,('BHR_Bahrain', 598)
    ,('BGR_Bulgaria', 598)
    ,('SDN_Sudan', 598)
    ,('POL_Poland', 598)
    ,('SAU_Saudi Arabia', 598)
    ,('IRQ_Iraq', 598)
    ,('HUN_Hungary', 598)
    ,('AGO_Angola', 598)
    ,('SWZ_Eswatini', 634)
    ,('MAC_China, Macao SAR', 634)
    ,('LAO_Lao People''s DR', 634)
    ,('MNG_Mongolia', 634)
    ,('BGD_Bangladesh', 654)
    ,('NPL_Nepal', 672)
    ,('MDG_Madagascar', 672)
    ,('HTI_Haiti', 672)
    ,('COG_Congo', 672)
    ,('SYR_Syrian Arab Republic', 672)
    ,('GMB_Gambia', 672)
    ,('DZA_Algeria', 672)
    ,('MLI_Mali', 672)
    ,('MMR_Myanmar', 708)
    ,('LBR_Liberia', 744)
    ,('GIN_Guinea', 826)
    ,('TCD_Chad', 841)
    ,('CPV_Cabo Verde', 841)
    ,('COM_Comoros', 841)
    ,('GNQ_Equatorial Guinea', 841)
    ,('GNB_Guinea-Bissau', 841)
    ,('SYC_Seychelles', 841)
    ,('VNM_Viet Nam', 852)
    ,('KHM_Cambodia', 852)
    ,('BLZ_Belize', 852)
    ,('ALB_Albania', 852)
    ,('MDV_Maldives', 852)
    ,('ARE_United Arab Emirates', 852)
    ,('GUY_Guyana', 852)
    ,('BRN_Brunei Darussalam', 852)
    ,('DJI_Djibouti', 991)
    ,('MSR_Montserrat', 991)
    ,('CYM_Cayman Islands', 991)
    ,('TCA_Turks and Caicos Islands', 991)
    ,('BTN_Bhutan', 991)
    ,('OMN_Oman', 991)
    ,('DMA_Dominica', 991)
    ,('BHS_Bahamas', 991)
    ,('ATG_Antigua and Barbuda', 991)
    ,('VCT_St. Vincent and the Grenadines', 991)
    ,('VGB_British Virgin Islands', 991)
    ,('AIA_Anguilla', 991)
    ,('BMU_Bermuda', 991)
    ,('SUR_Suriname', 991)
    ,('ABW_Aruba', 991)
    ,('LBN_Lebanon', 991)
    ,('GRD_Grenada', 991)
    ,('LCA_Saint Lucia', 991)
    ,('PSE_State of Palestine', 991)
    ,('STP_Sao Tome and Principe', 991)
    ,('KNA_Saint Kitts and Nevis', 991)
    ,('MDA_Republic of Moldova', 1078)
    ,('LVA_Latvia', 1078)
    ,('LTU_Lithuania', 1078)
    ,('ARM_Armenia', 1078)
    ,('UKR_Ukraine', 1078)
    ,('TJK_Tajikistan', 1078)
    ,('SRB_Serbia', 1078)
    ,('KGZ_Kyrgyzstan', 1078)
    ,('KAZ_Kazakhstan', 1078)
    ,('SVK_Slovakia', 1078)
    ,('HRV_Croatia', 1078)
    ,('CZE_Czech Republic', 1078)
    ,('RUS_Russian Federation', 1078)
    ,('SVN_Slovenia', 1078)
    ,('EST_Estonia', 1078)
    ,('YEM_Yemen', 1194)
    ,('TKM_Turkmenistan', 1291)
    ,('MNE_Montenegro', 1291)
    ,('MKD_North Macedonia', 1291)
    ,('GEO_Georgia', 1291)
    ,('AZE_Azerbaijan', 1291)
    ,('BIH_Bosnia and Herzegovina', 1291)
    ,('BLR_Belarus', 1291)
    ,('UZB_Uzbekistan', 1291)
    ,('CUW_Curaçao', 1516)
    ,('SXM_Sint Maarten (Dutch part)', 1516)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}